$d = $word.ActiveDocument

# Find the paragraph that holds "Ver no Jupiter Salvar em pdf Salvar em docx".
# That paragraph is immediately preceded by a blank paragraph (which also needs
# to go), and immediately followed by the "(c) 2020 ..." footer paragraph
# (which also needs to go). Deleting the range that spans from the start of
# the blank paragraph through the end of the footer paragraph collapses all
# three away, leaving the "LOB1255: ... (Requisito)" paragraph followed
# directly by the paragraph that used to come right after the footer.

$marker = $d.Content.Find
$marker.Text = "Ver no Jupiter Salvar em pdf Salvar em docx"
$marker.Execute() | Out-Null

$jupiterPara = $marker.Parent.Paragraphs(1)
$blankPara = $jupiterPara.Previous()
$footerPara = $jupiterPara.Next()

$start = $blankPara.Range.Start
$end = $footerPara.Range.End

$d.Range($start, $end).Delete()

$d.Save()
